$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E14").Value = "Framebuffer RenderAPI Employment"
$ws.Range("E14").Style = $ws.Range("B2").Style

$ws.Range("E14").Select()
